$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Keep the gridlines visible (runtime defaults to hidden on any rewrite
# unless this is explicitly restated) ---
$excel.ActiveWindow.DisplayGridlines = $true

# --- 1) Fix up the two bonito-flake strings that changed wording ---
$ws.Cells.Item(60, 3).Value2 = "かつおぶし変換を即座に有効化"
$ws.Cells.Item(60, 4).Value2 = "启用即时柴鱼片转换"
$ws.Cells.Item(61, 3).Value2 = "選択した魚をかつおぶしに即座に変換する機能を有効または無効にします。"

# --- 2) Column E was only ever used to carry a bare style on a few rows;
# clear it out so those cells (and the column) disappear from the sheet ---
$ws.Range("E14:E17").Clear()

# --- 3) Append the two new rows describing the "Exclude Tier Fish" toggle ---
$row1 = 103
$ws.Cells.Item($row1, 1).Value2 = "toggle53"
$ws.Cells.Item($row1, 2).Value2 = "Enable Exclude Tier Fish"

$c1jp = $ws.Cells.Item($row1, 3)
$c1jp.Value2 = "星付き魚を除外"

$c1cn = $ws.Cells.Item($row1, 4)
$c1cn.Value2 = "排除星级鱼类"

$row2 = 104
$ws.Cells.Item($row2, 1).Value2 = "tooltip19"
$ws.Cells.Item($row2, 2).Value2 = "Enable or disable excluding star-tier (★) fish from the bonito flakes feature.`nSet to 'true' to skip converting rare fish marked with a star."

$jpRun1 = "星付き（★）のレア魚をかつおぶし変換機能の対象から除外するかどうかを設定します。`n"
$jpRun2 = "'true' "
$jpRun3 = "に設定すると、★付きの魚は変換されません。"
$jpFull = $jpRun1 + $jpRun2 + $jpRun3

$c2jp = $ws.Cells.Item($row2, 3)
$c2jp.Value2 = $jpFull
$c2jp.Characters(1, $jpRun1.Length).Font.Name = "Noto Sans SC"
$c2jp.Characters(1, $jpRun1.Length).Font.Size = 10
$c2jp.Characters(1, $jpRun1.Length).Font.ColorIndex = -4105
$c2jp.Characters($jpRun1.Length + 1, $jpRun2.Length).Font.Name = "Arial"
$c2jp.Characters($jpRun1.Length + 1, $jpRun2.Length).Font.Size = 10
$c2jp.Characters($jpRun1.Length + 1, $jpRun2.Length).Font.ColorIndex = -4105
$c2jp.Characters($jpRun1.Length + $jpRun2.Length + 1, $jpRun3.Length).Font.Name = "Noto Sans SC"
$c2jp.Characters($jpRun1.Length + $jpRun2.Length + 1, $jpRun3.Length).Font.Size = 10
$c2jp.Characters($jpRun1.Length + $jpRun2.Length + 1, $jpRun3.Length).Font.ColorIndex = -4105

$cnRun1 = "启用或禁用将带有星标（★）的稀有鱼类排除在柴鱼片转换功能之外。`n设置为 "
$cnRun2 = "'true' "
$cnRun3 = "可跳过转换带星的稀有鱼类。"
$cnFull = $cnRun1 + $cnRun2 + $cnRun3

$c2cn = $ws.Cells.Item($row2, 4)
$c2cn.Value2 = $cnFull
$c2cn.Characters(1, $cnRun1.Length).Font.Name = "Noto Sans SC"
$c2cn.Characters(1, $cnRun1.Length).Font.Size = 10
$c2cn.Characters(1, $cnRun1.Length).Font.ColorIndex = -4105
$c2cn.Characters($cnRun1.Length + 1, $cnRun2.Length).Font.Name = "Arial"
$c2cn.Characters($cnRun1.Length + 1, $cnRun2.Length).Font.Size = 10
$c2cn.Characters($cnRun1.Length + 1, $cnRun2.Length).Font.ColorIndex = -4105
$c2cn.Characters($cnRun1.Length + $cnRun2.Length + 1, $cnRun3.Length).Font.Name = "Noto Sans SC"
$c2cn.Characters($cnRun1.Length + $cnRun2.Length + 1, $cnRun3.Length).Font.Size = 10
$c2cn.Characters($cnRun1.Length + $cnRun2.Length + 1, $cnRun3.Length).Font.ColorIndex = -4105

# Columns C/D on the new rows take the wrapped "Noto Sans SC" style.
$wrapRange = $ws.Range("C103:D104")
$wrapRange.Font.Name = "Noto Sans SC"
$wrapRange.Font.Size = 10
$wrapRange.WrapText = $true

# Row heights: single-line row keeps the sheet default, the wrapped
# tooltip row grows to fit its (roughly) three wrapped lines.
$ws.Rows.Item($row1).RowHeight = 12.8
$ws.Rows.Item($row2).RowHeight = 37.3

# --- 4) Restore the cursor/selection position shown in the diff ---
$ws.Range("D108").Select()
